$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F4 currently holds the inline string "111" -- convert it to a real number
$ws.Cells.Item(4, 6).Value = 111

# Add new row 5 with the additional user record
$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "dsadasd"
$ws.Cells.Item(5, 3).Value = "asdsadas"
$ws.Cells.Item(5, 4).Value = "dasdas"
$ws.Cells.Item(5, 5).Value = "chyrka"
$ws.Cells.Item(5, 6).Value = "dsadas"
$ws.Cells.Item(5, 7).Value = "dasdas"

# Match the style used by the other rows in column A (bold/border/centered)
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
